# Remove the "WolframAlpha Value" column (column D) contents from the
# ERF.PRECISE worksheet: clear the header in D1 and the values in D2:D22,
# then update the selection to match the post-edit state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "WolframAlpha Value" header cell.
$ws.Range("D1").ClearContents()

# Clear the WolframAlpha reference values themselves.
$ws.Range("D2:D22").ClearContents()

# Match the saved selection state (D1 active cell, D1:D22 selected).
$ws.Activate()
[void]$ws.Range("D1:D22").Select()
